$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 219 - shifts existing rows 219:265 down to 220:266
$ws.Rows.Item(219).Insert()

# Populate the newly inserted row 219 with the new weekly record.
# (Columns A,B,C,E,F,G,H,I,J,N,O,Q,R repeat the same "Ciboulette" series
# metadata as the row that used to sit at 219; D,K,L,M,P carry the new
# week's figures.)
$ws.Cells.Item(219, 1).Value = 9
$ws.Cells.Item(219, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(219, 3).Value = "Metropolitana"
$ws.Cells.Item(219, 4).Value = 44476
$ws.Cells.Item(219, 5).Value = 13
$ws.Cells.Item(219, 6).Value = 100112039
$ws.Cells.Item(219, 7).Value = "Ciboulette"
$ws.Cells.Item(219, 8).Value = "Sin especificar"
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 250
$ws.Cells.Item(219, 11).Value = 800
$ws.Cells.Item(219, 12).Value = 1000
$ws.Cells.Item(219, 13).Value = 900
$ws.Cells.Item(219, 14).Value = "`$/docena de atados"
$ws.Cells.Item(219, 15).Value = "Región Metropolitana"
$ws.Cells.Item(219, 16).Value = 300
$ws.Cells.Item(219, 17).Value = 3
$ws.Cells.Item(219, 18).Value = "Hortaliza"
